$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new blank "Strings" sheet right after Sheet1 first (so it gets the
# lower internal sheetId), then copy Sheet1 to produce "Sheet1 (2)" which will
# be inserted between Sheet1 and Strings.
$stringsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$stringsSheet.Name = "Strings"
[void]$stringsSheet.Range("A2").Select()

$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Update the content that differs on the copied sheet (order matches the
# shared-string table order: Strings, Longest Palindromic Substring, Doubt,
# Do it by DP way or by interview camp method).
$ws2.Range("B3").Value = "Strings"
$ws2.Range("C3").Value = "Longest Palindromic Substring"
$ws2.Range("D1").Value = "Doubt"
$ws2.Range("D3").Value = "Do it by DP way or by interview camp method"

# Restore selection on the original sheet and move it off of the active tab.
[void]$ws1.Range("A14").Select()

# Make "Sheet1 (2)" the active sheet with its own selection.
[void]$ws2.Activate()
[void]$ws2.Range("C5").Select()
